# Update cryptos list (price / 1h volume change) and swap the
# ApeXProtocol / Stellar rows (47 and 48) to reflect the new ranking.
# NumberFormat is forced to "@" (Text) before each write so that
# numeric-looking price strings (e.g. "0.614", "12.42") are stored as
# literal text, matching the original inlineStr cell content instead
# of being auto-converted to floating point numbers by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.763.97"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.27%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.564.75"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -4.22%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.64"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -6.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "187.25"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.561.71"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -4.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.614"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.86%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.671"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -7.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.149"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -8.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.52"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -7.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000263"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -9.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.80"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -6.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.133.48"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.558.63"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -4.35%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.28"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -6.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "66.712.92"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.09"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -6.85%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -8.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "390.14"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -5.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.23"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -7.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.55"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -5.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.21"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -5.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.42"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -5.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.07"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.62%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -7.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.87"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -9.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.56"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "30.92"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -6.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "628.25"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.19"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.65%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -7.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "63.33"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -6.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "41.84"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -9.47%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.58%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0757"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -9.70%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.102.06"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +6.16%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.93"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.65"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0410"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -8.53%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.11"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.86%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.130"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -7.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "139.01"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.78"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.47"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -8.95%  "
